# Update the LR-pairs data (Pf4-Fgfr2) sheet with newly recomputed TPM-based values.
# The Target cluster grouping changed from Inflammatory-Mac to Resolving-Mac for one
# quarter of rows, and all dependent statistics (columns E:T) were recalculated
# against the new TPM data, as reflected cell-by-cell below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'Pf4'
$ws.Cells.Item(2,3).Value = 'Fgfr2'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 5.090240666666667
$ws.Cells.Item(2,8).Value = 15.270722
$ws.Cells.Item(2,9).Value = 0.01518526656315472
$ws.Cells.Item(2,10).Value = 0.01525191836740238
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.2858606666666667
$ws.Cells.Item(2,14).Value = 0.857582
$ws.Cells.Item(2,15).Value = 0.0687156860066334
$ws.Cells.Item(2,16).Value = 0.06932858672617494
$ws.Cells.Item(2,17).Value = 1.455099590467111
$ws.Cells.Item(2,18).Value = 13.095896314204
$ws.Cells.Item(2,19).Value = 0.001043466009080769
$ws.Cells.Item(2,20).Value = 0.001057393945274997
$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'Pf4'
$ws.Cells.Item(3,3).Value = 'Fgfr2'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 5.090240666666667
$ws.Cells.Item(3,8).Value = 15.270722
$ws.Cells.Item(3,9).Value = 0.01518526656315472
$ws.Cells.Item(3,10).Value = 0.01525191836740238
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.763360333333333
$ws.Cells.Item(3,14).Value = 11.290081
$ws.Cells.Item(3,15).Value = 0.90464312565499
$ws.Cells.Item(3,16).Value = 0.9127119736118995
$ws.Cells.Item(3,17).Value = 19.15640981205356
$ws.Cells.Item(3,18).Value = 172.407688308482
$ws.Cells.Item(3,19).Value = 0.01373724700759649
$ws.Cells.Item(3,20).Value = 0.01392060851447941
$ws.Cells.Item(4,1).Value = 'ECs'
$ws.Cells.Item(4,2).Value = 'Pf4'
$ws.Cells.Item(4,3).Value = 'Fgfr2'
$ws.Cells.Item(4,4).Value = 'MuSCs'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 5.090240666666667
$ws.Cells.Item(4,8).Value = 15.270722
$ws.Cells.Item(4,9).Value = 0.01518526656315472
$ws.Cells.Item(4,10).Value = 0.01525191836740238
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.110331
$ws.Cells.Item(4,14).Value = 0.220662
$ws.Cells.Item(4,15).Value = 0.02652155835639462
$ws.Cells.Item(4,16).Value = 0.01783874265571248
$ws.Cells.Item(4,17).Value = 0.5616113429940001
$ws.Cells.Item(4,18).Value = 3.369668057964
$ws.Cells.Item(4,19).Value = 0.0004027369333121158
$ws.Cells.Item(4,20).Value = 0.0002720750467620254
$ws.Cells.Item(5,1).Value = 'ECs'
$ws.Cells.Item(5,2).Value = 'Pf4'
$ws.Cells.Item(5,3).Value = 'Fgfr2'
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 5.090240666666667
$ws.Cells.Item(5,8).Value = 15.270722
$ws.Cells.Item(5,9).Value = 0.01518526656315472
$ws.Cells.Item(5,10).Value = 0.01525191836740238
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.0004976666666666667
$ws.Cells.Item(5,14).Value = 0.001493
$ws.Cells.Item(5,15).Value = 0.0001196299819817856
$ws.Cells.Item(5,16).Value = 0.0001206970062130259
$ws.Cells.Item(5,17).Value = 0.002533243105111111
$ws.Cells.Item(5,18).Value = 0.022799187946
$ws.Cells.Item(5,19).Value = [double]"1.816613165338811E-06"
$ws.Cells.Item(5,20).Value = [double]"1.840860885950929E-06"
$ws.Cells.Item(6,1).Value = 'FAPs'
$ws.Cells.Item(6,2).Value = 'Pf4'
$ws.Cells.Item(6,3).Value = 'Fgfr2'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.020288666666667
$ws.Cells.Item(6,8).Value = 3.060866
$ws.Cells.Item(6,9).Value = 0.003043737298347591
$ws.Cells.Item(6,10).Value = 0.003057096996825524
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.2858606666666667
$ws.Cells.Item(6,14).Value = 0.857582
$ws.Cells.Item(6,15).Value = 0.0687156860066334
$ws.Cells.Item(6,16).Value = 0.06932858672617494
$ws.Cells.Item(6,17).Value = 0.2916603984457777
$ws.Cells.Item(6,18).Value = 2.624943586012
$ws.Cells.Item(6,19).Value = 0.0002091524964799317
$ws.Cells.Item(6,20).Value = 0.0002119442142747473
$ws.Cells.Item(7,1).Value = 'FAPs'
$ws.Cells.Item(7,2).Value = 'Pf4'
$ws.Cells.Item(7,3).Value = 'Fgfr2'
$ws.Cells.Item(7,4).Value = 'FAPs'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.020288666666667
$ws.Cells.Item(7,8).Value = 3.060866
$ws.Cells.Item(7,9).Value = 0.003043737298347591
$ws.Cells.Item(7,10).Value = 0.003057096996825524
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.763360333333333
$ws.Cells.Item(7,14).Value = 11.290081
$ws.Cells.Item(7,15).Value = 0.90464312565499
$ws.Cells.Item(7,16).Value = 0.9127119736118995
$ws.Cells.Item(7,17).Value = 3.839713896682889
$ws.Cells.Item(7,18).Value = 34.557425070146
$ws.Cells.Item(7,19).Value = 0.00275349602324984
$ws.Cells.Item(7,20).Value = 0.002790249033495635
$ws.Cells.Item(8,1).Value = 'FAPs'
$ws.Cells.Item(8,2).Value = 'Pf4'
$ws.Cells.Item(8,3).Value = 'Fgfr2'
$ws.Cells.Item(8,4).Value = 'MuSCs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.020288666666667
$ws.Cells.Item(8,8).Value = 3.060866
$ws.Cells.Item(8,9).Value = 0.003043737298347591
$ws.Cells.Item(8,10).Value = 0.003057096996825524
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.110331
$ws.Cells.Item(8,14).Value = 0.220662
$ws.Cells.Item(8,15).Value = 0.02652155835639462
$ws.Cells.Item(8,16).Value = 0.01783874265571248
$ws.Cells.Item(8,17).Value = 0.112569468882
$ws.Cells.Item(8,18).Value = 0.675416813292
$ws.Cells.Item(8,19).Value = [double]"8.072465637966054E-05"
$ws.Cells.Item(8,20).Value = [double]"5.453476659992197E-05"
$ws.Cells.Item(9,1).Value = 'FAPs'
$ws.Cells.Item(9,2).Value = 'Pf4'
$ws.Cells.Item(9,3).Value = 'Fgfr2'
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.020288666666667
$ws.Cells.Item(9,8).Value = 3.060866
$ws.Cells.Item(9,9).Value = 0.003043737298347591
$ws.Cells.Item(9,10).Value = 0.003057096996825524
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.0004976666666666667
$ws.Cells.Item(9,14).Value = 0.001493
$ws.Cells.Item(9,15).Value = 0.0001196299819817856
$ws.Cells.Item(9,16).Value = 0.0001206970062130259
$ws.Cells.Item(9,17).Value = 0.0005077636597777777
$ws.Cells.Item(9,18).Value = 0.004569872938
$ws.Cells.Item(9,19).Value = [double]"3.641222381586111E-07"
$ws.Cells.Item(9,20).Value = [double]"3.68982455219673E-07"
$ws.Cells.Item(10,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(10,2).Value = 'Pf4'
$ws.Cells.Item(10,3).Value = 'Fgfr2'
$ws.Cells.Item(10,4).Value = 'ECs'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 69.13821133333333
$ws.Cells.Item(10,8).Value = 207.414634
$ws.Cells.Item(10,9).Value = 0.206253935235621
$ws.Cells.Item(10,10).Value = 0.2071592335956769
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.2858606666666667
$ws.Cells.Item(10,14).Value = 0.857582
$ws.Cells.Item(10,15).Value = 0.0687156860066334
$ws.Cells.Item(10,16).Value = 0.06932858672617494
$ws.Cells.Item(10,17).Value = 19.76389518388755
$ws.Cells.Item(10,18).Value = 177.875056654988
$ws.Cells.Item(10,19).Value = 0.01417288065128344
$ws.Cells.Item(10,20).Value = 0.01436205689246582
$ws.Cells.Item(11,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(11,2).Value = 'Pf4'
$ws.Cells.Item(11,3).Value = 'Fgfr2'
$ws.Cells.Item(11,4).Value = 'FAPs'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 69.13821133333333
$ws.Cells.Item(11,8).Value = 207.414634
$ws.Cells.Item(11,9).Value = 0.206253935235621
$ws.Cells.Item(11,10).Value = 0.2071592335956769
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 3.763360333333333
$ws.Cells.Item(11,14).Value = 11.290081
$ws.Cells.Item(11,15).Value = 0.90464312565499
$ws.Cells.Item(11,16).Value = 0.9127119736118995
$ws.Cells.Item(11,17).Value = 260.1920020494838
$ws.Cells.Item(11,18).Value = 2341.728018445354
$ws.Cells.Item(11,19).Value = 0.1865862046501941
$ws.Cells.Item(11,20).Value = 0.1890767129470388
$ws.Cells.Item(12,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(12,2).Value = 'Pf4'
$ws.Cells.Item(12,3).Value = 'Fgfr2'
$ws.Cells.Item(12,4).Value = 'MuSCs'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 69.13821133333333
$ws.Cells.Item(12,8).Value = 207.414634
$ws.Cells.Item(12,9).Value = 0.206253935235621
$ws.Cells.Item(12,10).Value = 0.2071592335956769
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.110331
$ws.Cells.Item(12,14).Value = 0.220662
$ws.Cells.Item(12,15).Value = 0.02652155835639462
$ws.Cells.Item(12,16).Value = 0.01783874265571248
$ws.Cells.Item(12,17).Value = 7.628087994617999
$ws.Cells.Item(12,18).Value = 45.768527967708
$ws.Cells.Item(12,19).Value = 0.005470175779587559
$ws.Cells.Item(12,20).Value = 0.003695460256867906
$ws.Cells.Item(13,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(13,2).Value = 'Pf4'
$ws.Cells.Item(13,3).Value = 'Fgfr2'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 69.13821133333333
$ws.Cells.Item(13,8).Value = 207.414634
$ws.Cells.Item(13,9).Value = 0.206253935235621
$ws.Cells.Item(13,10).Value = 0.2071592335956769
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.0004976666666666667
$ws.Cells.Item(13,14).Value = 0.001493
$ws.Cells.Item(13,15).Value = 0.0001196299819817856
$ws.Cells.Item(13,16).Value = 0.0001206970062130259
$ws.Cells.Item(13,17).Value = 0.03440778317355556
$ws.Cells.Item(13,18).Value = 0.309670048562
$ws.Cells.Item(13,19).Value = [double]"2.467415455590972E-05"
$ws.Cells.Item(13,20).Value = [double]"2.50034993043831E-05"
$ws.Cells.Item(14,1).Value = 'MuSCs'
$ws.Cells.Item(14,2).Value = 'Pf4'
$ws.Cells.Item(14,3).Value = 'Fgfr2'
$ws.Cells.Item(14,4).Value = 'ECs'
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.394653
$ws.Cells.Item(14,8).Value = 8.789306
$ws.Cells.Item(14,9).Value = 0.01311018115402158
$ws.Cells.Item(14,10).Value = 0.008778483271329277
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.2858606666666667
$ws.Cells.Item(14,14).Value = 0.857582
$ws.Cells.Item(14,15).Value = 0.0687156860066334
$ws.Cells.Item(14,16).Value = 0.06932858672617494
$ws.Cells.Item(14,17).Value = 1.256258436348667
$ws.Cells.Item(14,18).Value = 7.537550618091999
$ws.Cells.Item(14,19).Value = 0.0009008750916698295
$ws.Cells.Item(14,20).Value = 0.0006085998388006276
$ws.Cells.Item(15,1).Value = 'MuSCs'
$ws.Cells.Item(15,2).Value = 'Pf4'
$ws.Cells.Item(15,3).Value = 'Fgfr2'
$ws.Cells.Item(15,4).Value = 'FAPs'
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.394653
$ws.Cells.Item(15,8).Value = 8.789306
$ws.Cells.Item(15,9).Value = 0.01311018115402158
$ws.Cells.Item(15,10).Value = 0.008778483271329277
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.763360333333333
$ws.Cells.Item(15,14).Value = 11.290081
$ws.Cells.Item(15,15).Value = 0.90464312565499
$ws.Cells.Item(15,16).Value = 0.9127119736118995
$ws.Cells.Item(15,17).Value = 16.53866277896433
$ws.Cells.Item(15,18).Value = 99.23197667378601
$ws.Cells.Item(15,19).Value = 0.01186003525707723
$ws.Cells.Item(15,20).Value = 0.008012226791893989
$ws.Cells.Item(16,1).Value = 'MuSCs'
$ws.Cells.Item(16,2).Value = 'Pf4'
$ws.Cells.Item(16,3).Value = 'Fgfr2'
$ws.Cells.Item(16,4).Value = 'MuSCs'
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.394653
$ws.Cells.Item(16,8).Value = 8.789306
$ws.Cells.Item(16,9).Value = 0.01311018115402158
$ws.Cells.Item(16,10).Value = 0.008778483271329277
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.110331
$ws.Cells.Item(16,14).Value = 0.220662
$ws.Cells.Item(16,15).Value = 0.02652155835639462
$ws.Cells.Item(16,16).Value = 0.01783874265571248
$ws.Cells.Item(16,17).Value = 0.484866460143
$ws.Cells.Item(16,18).Value = 1.939465840572
$ws.Cells.Item(16,19).Value = 0.0003477024345392882
$ws.Cells.Item(16,20).Value = 0.00015659710398472
$ws.Cells.Item(17,1).Value = 'MuSCs'
$ws.Cells.Item(17,2).Value = 'Pf4'
$ws.Cells.Item(17,3).Value = 'Fgfr2'
$ws.Cells.Item(17,4).Value = 'Resolving-Mac'
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.394653
$ws.Cells.Item(17,8).Value = 8.789306
$ws.Cells.Item(17,9).Value = 0.01311018115402158
$ws.Cells.Item(17,10).Value = 0.008778483271329277
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.0004976666666666667
$ws.Cells.Item(17,14).Value = 0.001493
$ws.Cells.Item(17,15).Value = 0.0001196299819817856
$ws.Cells.Item(17,16).Value = 0.0001206970062130259
$ws.Cells.Item(17,17).Value = 0.002187072309666666
$ws.Cells.Item(17,18).Value = 0.013122433858
$ws.Cells.Item(17,19).Value = [double]"1.568370735233547E-06"
$ws.Cells.Item(17,20).Value = [double]"1.059536649940574E-06"
$ws.Cells.Item(18,1).Value = 'Resolving-Mac'
$ws.Cells.Item(18,2).Value = 'Pf4'
$ws.Cells.Item(18,3).Value = 'Fgfr2'
$ws.Cells.Item(18,4).Value = 'ECs'
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 255.5657806666667
$ws.Cells.Item(18,8).Value = 766.6973419999999
$ws.Cells.Item(18,9).Value = 0.7624068797488551
$ws.Cells.Item(18,10).Value = 0.765753267768766
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 0.2858606666666667
$ws.Cells.Item(18,14).Value = 0.857582
$ws.Cells.Item(18,15).Value = 0.0687156860066334
$ws.Cells.Item(18,16).Value = 0.06932858672617494
$ws.Cells.Item(18,17).Value = 73.05620443856044
$ws.Cells.Item(18,18).Value = 657.5058399470439
$ws.Cells.Item(18,19).Value = 0.05238931175811944
$ws.Cells.Item(18,20).Value = 0.05308859183535875
$ws.Cells.Item(19,1).Value = 'Resolving-Mac'
$ws.Cells.Item(19,2).Value = 'Pf4'
$ws.Cells.Item(19,3).Value = 'Fgfr2'
$ws.Cells.Item(19,4).Value = 'FAPs'
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 255.5657806666667
$ws.Cells.Item(19,8).Value = 766.6973419999999
$ws.Cells.Item(19,9).Value = 0.7624068797488551
$ws.Cells.Item(19,10).Value = 0.765753267768766
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 3.763360333333333
$ws.Cells.Item(19,14).Value = 11.290081
$ws.Cells.Item(19,15).Value = 0.90464312565499
$ws.Cells.Item(19,16).Value = 0.9127119736118995
$ws.Cells.Item(19,17).Value = 961.7861215183002
$ws.Cells.Item(19,18).Value = 8656.075093664702
$ws.Cells.Item(19,19).Value = 0.6897061427168725
$ws.Cells.Item(19,20).Value = 0.6989121763249917
$ws.Cells.Item(20,1).Value = 'Resolving-Mac'
$ws.Cells.Item(20,2).Value = 'Pf4'
$ws.Cells.Item(20,3).Value = 'Fgfr2'
$ws.Cells.Item(20,4).Value = 'MuSCs'
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 255.5657806666667
$ws.Cells.Item(20,8).Value = 766.6973419999999
$ws.Cells.Item(20,9).Value = 0.7624068797488551
$ws.Cells.Item(20,10).Value = 0.765753267768766
$ws.Cells.Item(20,11).Value = 2
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 0.110331
$ws.Cells.Item(20,14).Value = 0.220662
$ws.Cells.Item(20,15).Value = 0.02652155835639462
$ws.Cells.Item(20,16).Value = 0.01783874265571248
$ws.Cells.Item(20,17).Value = 28.196828146734
$ws.Cells.Item(20,18).Value = 169.180968880404
$ws.Cells.Item(20,19).Value = 0.02022021855257599
$ws.Cells.Item(20,20).Value = 0.0136600754814979
$ws.Cells.Item(21,1).Value = 'Resolving-Mac'
$ws.Cells.Item(21,2).Value = 'Pf4'
$ws.Cells.Item(21,3).Value = 'Fgfr2'
$ws.Cells.Item(21,4).Value = 'Resolving-Mac'
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 255.5657806666667
$ws.Cells.Item(21,8).Value = 766.6973419999999
$ws.Cells.Item(21,9).Value = 0.7624068797488551
$ws.Cells.Item(21,10).Value = 0.765753267768766
$ws.Cells.Item(21,11).Value = 1
$ws.Cells.Item(21,12).Value = 0.3333333333333333
$ws.Cells.Item(21,13).Value = 0.0004976666666666667
$ws.Cells.Item(21,14).Value = 0.001493
$ws.Cells.Item(21,15).Value = 0.0001196299819817856
$ws.Cells.Item(21,16).Value = 0.0001206970062130259
$ws.Cells.Item(21,17).Value = 0.1271865701784444
$ws.Cells.Item(21,18).Value = 1.144679131606
$ws.Cells.Item(21,19).Value = [double]"9.120672128714493E-05"
$ws.Cells.Item(21,20).Value = [double]"9.242412691753163E-05"
